{"js": "const replacements = [\n  [\"403\u00f73=134, 1\", \"875\u00f77=125, 0\"],\n  [\"599\u00f78=74, 7\", \"617\u00f77=88, 1\"],\n  [\"275\u00f72=137, 1\", \"208\u00f78=26, 0\"],\n  [\"269\u00f79=29, 8\", \"605\u00f74=151, 1\"],\n  [\"541\u00f76=90, 1\", \"105\u00f77=15, 0\"],\n  [\"795\u00f78=99, 3\", \"732\u00f78=91, 4\"],\n  [\"545\u00f75=109, 0\", \"824\u00f73=274, 2\"],\n  [\"405\u00f72=202, 1\", \"998\u00f79=110, 8\"],\n  [\"783\u00f78=97, 7\", \"627\u00f75=125, 2\"],\n  [\"393\u00f74=98, 1\", \"639\u00f75=127, 4\"],\n  [\"880\u00f79=97, 7\", \"743\u00f77=106, 1\"],\n  [\"346\u00f74=86, 2\", \"671\u00f76=111, 5\"],\n  [\"879\u00f78=109, 7\", \"112\u00f72=56, 0\"],\n  [\"439\u00f73=146, 1\", \"562\u00f76=93, 4\"],\n  [\"937\u00f73=312, 1\", \"399\u00f78=49, 7\"],\n  [\"341\u00f73=113, 2\", \"114\u00f73=38, 0\"],\n  [\"897\u00f77=128, 1\", \"402\u00f73=134, 0\"],\n  [\"580\u00f75=116, 0\", \"204\u00f74=51, 0\"],\n  [\"451\u00f74=112, 3\", \"769\u00f74=192, 1\"],\n  [\"834\u00f76=139, 0\", \"657\u00f79=73, 0\"],\n  [\"516\u00f77=73, 5\", \"981\u00f78=122, 5\"],\n  [\"219\u00f73=73, 0\", \"306\u00f77=43, 5\"],\n  [\"936\u00f76=156, 0\", \"593\u00f75=118, 3\"],\n  [\"968\u00f76=161, 2\", \"698\u00f79=77, 5\"],\n  [\"504\u00f78=63, 0\", \"652\u00f74=163, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"403\u00f73=134, 1\", \"875\u00f77=125, 0\"),\n    @(\"599\u00f78=74, 7\", \"617\u00f77=88, 1\"),\n    @(\"275\u00f72=137, 1\", \"208\u00f78=26, 0\"),\n    @(\"269\u00f79=29, 8\", \"605\u00f74=151, 1\"),\n    @(\"541\u00f76=90, 1\", \"105\u00f77=15, 0\"),\n    @(\"795\u00f78=99, 3\", \"732\u00f78=91, 4\"),\n    @(\"545\u00f75=109, 0\", \"824\u00f73=274, 2\"),\n    @(\"405\u00f72=202, 1\", \"998\u00f79=110, 8\"),\n    @(\"783\u00f78=97, 7\", \"627\u00f75=125, 2\"),\n    @(\"393\u00f74=98, 1\", \"639\u00f75=127, 4\"),\n    @(\"880\u00f79=97, 7\", \"743\u00f77=106, 1\"),\n    @(\"346\u00f74=86, 2\", \"671\u00f76=111, 5\"),\n    @(\"879\u00f78=109, 7\", \"112\u00f72=56, 0\"),\n    @(\"439\u00f73=146, 1\", \"562\u00f76=93, 4\"),\n    @(\"937\u00f73=312, 1\", \"399\u00f78=49, 7\"),\n    @(\"341\u00f73=113, 2\", \"114\u00f73=38, 0\"),\n    @(\"897\u00f77=128, 1\", \"402\u00f73=134, 0\"),\n    @(\"580\u00f75=116, 0\", \"204\u00f74=51, 0\"),\n    @(\"451\u00f74=112, 3\", \"769\u00f74=192, 1\"),\n    @(\"834\u00f76=139, 0\", \"657\u00f79=73, 0\"),\n    @(\"516\u00f77=73, 5\", \"981\u00f78=122, 5\"),\n    @(\"219\u00f73=73, 0\", \"306\u00f77=43, 5\"),\n    @(\"936\u00f76=156, 0\", \"593\u00f75=118, 3\"),\n    @(\"968\u00f76=161, 2\", \"698\u00f79=77, 5\"),\n    @(\"504\u00f78=63, 0\", \"652\u00f74=163, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
